$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '329.18'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.30%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.30'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.15%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.51%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08098'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.23%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.981'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '4.69%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.320'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9525'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.37%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1178'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-3.91%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1850'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.55%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'MCDex'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '10.24'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '19.98%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09809'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.47%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04633'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '12.29%'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.1068'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.11%'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001285'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.93%'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04219'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.94%'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005873'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.93%'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.372'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-5.59%'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3472'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.74%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '5.92%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2507'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.51%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.11%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004319'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.42%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.60%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02656'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '0.29%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05558'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.75%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007567'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.56%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1409'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.008087'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-16.80%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.82%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008887'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-10.28%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007236'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.67%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000751'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.27%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.002272'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.27%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '21.88%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002102'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.27%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002002'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.27%'
